# Change the "Runmode" values for the Login functionality and ProductListPage
# test rows from "Y" to "N" (object for where-to-buy TC changed), and move
# the active selection to B7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "N"
$ws.Range("C3").Value = "N"

$ws.Range("B7").Select()
